# Renaming of ecology_format tables
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ecological_params")
$ws2 = $wb.Worksheets.Item("Scaling")

# Rename the PLA cuboid headers on the ecological_params sheet
$hdr = $ws1.Range("B1:D1")
$ws1.Range("B1").Value = "PLA_virgin"
$ws1.Range("C1").Value = "PLA_recycled"
$ws1.Range("D1").Value = "PLA_recycled_industrial"

# Re-style the renamed header cells: plain Times New Roman, no border / special alignment
$hdr.Style = "Normal"
$hdr.Font.Name = "Times New Roman"
$hdr.Font.Size = 12
$hdr.Font.Color = 0

# ecological_params becomes the active / selected sheet, with C8 selected
$ws1.Activate()
$ws1.Range("C8").Select()
